$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("To Do")

# The "File not showing in the submitted list after adding" task is now
# also just "Done" (was the one-off "Already Done" status).
$ws.Range("B5").Value = "Done"

# The remaining admin-web tasks (rows 16-19) are now marked "Done" too -
# they previously had no status in column B.
$ws.Range("B16").Value = "Done"
$ws.Range("B17").Value = "Done"
$ws.Range("B18").Value = "Done"
$ws.Range("B19").Value = "Done"

# Update the view: zoom-to-100 indicator + move the active selection.
$ws.Application.ActiveWindow.Zoom = 100
$ws.Range("A18").Select()
